# cleartrip_testdata.xlsx - add the "restServiceForCleartrip" REST/JSON test-data sheet
# (SonarLint code-review follow-up: new worksheet with one-way flight search payload data)

$wb = $excel.ActiveWorkbook

# --- add the new worksheet after the existing two, make it the active sheet ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "restServiceForCleartrip"

# --- header row ---
$ws3.Range("A1").Value = "trip_type"
$ws3.Range("B1").Value = "source"
$ws3.Range("C1").Value = "destination"
$ws3.Range("D1").Value = "depart_date"
$ws3.Range("E1").Value = "adults"
$ws3.Range("F1").Value = "children"
$ws3.Range("G1").Value = "infants"
$ws3.Range("H1").Value = "origin"
$ws3.Range("I1").Value = "from"
$ws3.Range("J1").Value = "to"
$ws3.Range("L1").Value = "ver"
$ws3.Range("M1").Value = "type"

# --- data row ---
$ws3.Range("A2").Value = "ONEWAY"
$ws3.Range("B2").Value = "Pune"
$ws3.Range("C2").Value = "Mumbai IN Chatrapati Shivaji Airport (BOM)"
$ws3.Range("D2").Value = "25/10/2018"
$ws3.Range("E2").Value = 1
$ws3.Range("F2").Value = 1
$ws3.Range("G2").Value = 1
$ws3.Range("H2").Value = "Pune%2C+IN+-+Lohegaon+(PNQ)"
$ws3.Range("I2").Value = "PNQ"
$ws3.Range("J2").Value = "BOM"
$ws3.Range("K2").Value = "Economy"
$ws3.Range("L2").Value = "V2"
$ws3.Range("M2").Value = "JSON"

# K1 header entered last (matches authored shared-string ordering)
$ws3.Range("K1").Value = "class1"

# --- column widths (best-fit sizing for the widest entries in each column) ---
$ws3.Columns("C:C").ColumnWidth = 39.67
$ws3.Columns("D:D").ColumnWidth = 11
$ws3.Columns("E:E").ColumnWidth = 5.67
$ws3.Columns("F:F").ColumnWidth = 7.5
$ws3.Columns("H:H").ColumnWidth = 29.83

# --- selection / active cell on the new sheet ---
$ws3.Range("J10").Select() | Out-Null

Write-Host "Added restServiceForCleartrip with" $wb.Worksheets.Count "sheets total"
